$d = $word.ActiveDocument

$replacements = @(
    @{old = "790×7=5530"; new = "865×6=5190"},
    @{old = "910×2=1820"; new = "253×9=2277"},
    @{old = "192×9=1728"; new = "827×9=7443"},
    @{old = "250×3=750";  new = "391×5=1955"},
    @{old = "754×7=5278"; new = "742×7=5194"},
    @{old = "216×8=1728"; new = "533×3=1599"},
    @{old = "332×6=1992"; new = "804×8=6432"},
    @{old = "404×5=2020"; new = "946×7=6622"},
    @{old = "248×6=1488"; new = "357×2=714"},
    @{old = "678×4=2712"; new = "931×5=4655"},
    @{old = "851×3=2553"; new = "328×8=2624"},
    @{old = "601×5=3005"; new = "398×4=1592"},
    @{old = "713×5=3565"; new = "845×8=6760"},
    @{old = "382×4=1528"; new = "730×2=1460"},
    @{old = "956×9=8604"; new = "251×8=2008"},
    @{old = "196×6=1176"; new = "890×8=7120"},
    @{old = "702×2=1404"; new = "357×3=1071"},
    @{old = "887×8=7096"; new = "204×9=1836"},
    @{old = "524×2=1048"; new = "746×2=1492"},
    @{old = "453×8=3624"; new = "590×5=2950"},
    @{old = "316×3=948";  new = "813×3=2439"},
    @{old = "616×6=3696"; new = "634×6=3804"},
    @{old = "652×3=1956"; new = "508×7=3556"},
    @{old = "606×9=5454"; new = "626×6=3756"},
    @{old = "654×6=3924"; new = "469×4=1876"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2) | Out-Null
}
